{"js": "const replacements = [\n  [\"762\u00d72=1524\", \"505\u00d72=1010\"],\n  [\"208\u00d73=624\", \"134\u00d78=1072\"],\n  [\"805\u00d72=1610\", \"571\u00d79=5139\"],\n  [\"872\u00d74=3488\", \"475\u00d73=1425\"],\n  [\"112\u00d73=336\", \"434\u00d79=3906\"],\n  [\"295\u00d78=2360\", \"573\u00d78=4584\"],\n  [\"921\u00d76=5526\", \"460\u00d75=2300\"],\n  [\"995\u00d76=5970\", \"368\u00d75=1840\"],\n  [\"422\u00d74=1688\", \"456\u00d74=1824\"],\n  [\"881\u00d72=1762\", \"293\u00d79=2637\"],\n  [\"909\u00d79=8181\", \"718\u00d72=1436\"],\n  [\"737\u00d73=2211\", \"293\u00d76=1758\"],\n  [\"949\u00d79=8541\", \"864\u00d73=2592\"],\n  [\"633\u00d79=5697\", \"380\u00d78=3040\"],\n  [\"711\u00d77=4977\", \"471\u00d76=2826\"],\n  [\"771\u00d77=5397\", \"955\u00d79=8595\"],\n  [\"848\u00d78=6784\", \"839\u00d79=7551\"],\n  [\"552\u00d75=2760\", \"566\u00d76=3396\"],\n  [\"183\u00d77=1281\", \"625\u00d73=1875\"],\n  [\"737\u00d79=6633\", \"620\u00d74=2480\"],\n  [\"660\u00d78=5280\", \"872\u00d78=6976\"],\n  [\"840\u00d78=6720\", \"623\u00d76=3738\"],\n  [\"349\u00d79=3141\", \"896\u00d74=3584\"],\n  [\"340\u00d75=1700\", \"892\u00d77=6244\"],\n  [\"126\u00d79=1134\", \"558\u00d74=2232\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"762\u00d72=1524\", \"505\u00d72=1010\"),\n    @(\"208\u00d73=624\", \"134\u00d78=1072\"),\n    @(\"805\u00d72=1610\", \"571\u00d79=5139\"),\n    @(\"872\u00d74=3488\", \"475\u00d73=1425\"),\n    @(\"112\u00d73=336\", \"434\u00d79=3906\"),\n    @(\"295\u00d78=2360\", \"573\u00d78=4584\"),\n    @(\"921\u00d76=5526\", \"460\u00d75=2300\"),\n    @(\"995\u00d76=5970\", \"368\u00d75=1840\"),\n    @(\"422\u00d74=1688\", \"456\u00d74=1824\"),\n    @(\"881\u00d72=1762\", \"293\u00d79=2637\"),\n    @(\"909\u00d79=8181\", \"718\u00d72=1436\"),\n    @(\"737\u00d73=2211\", \"293\u00d76=1758\"),\n    @(\"949\u00d79=8541\", \"864\u00d73=2592\"),\n    @(\"633\u00d79=5697\", \"380\u00d78=3040\"),\n    @(\"711\u00d77=4977\", \"471\u00d76=2826\"),\n    @(\"771\u00d77=5397\", \"955\u00d79=8595\"),\n    @(\"848\u00d78=6784\", \"839\u00d79=7551\"),\n    @(\"552\u00d75=2760\", \"566\u00d76=3396\"),\n    @(\"183\u00d77=1281\", \"625\u00d73=1875\"),\n    @(\"737\u00d79=6633\", \"620\u00d74=2480\"),\n    @(\"660\u00d78=5280\", \"872\u00d78=6976\"),\n    @(\"840\u00d78=6720\", \"623\u00d76=3738\"),\n    @(\"349\u00d79=3141\", \"896\u00d74=3584\"),\n    @(\"340\u00d75=1700\", \"892\u00d77=6244\"),\n    @(\"126\u00d79=1134\", \"558\u00d74=2232\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
